$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 909208.5600000001
$ws.Range("I5").Value = 1250118
$ws.Range("J5").Value = 116.666664
$ws.Range("K5").Value = 1250118
$ws.Range("L5").Value = 116.666664
$ws.Range("M5").Value = -1250003
$ws.Range("N5").Value = -346.666664

$ws.Range("H69").Value = 254711.25
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 254711.25
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 764133.75
$ws.Range("N69").Value = -765881.75
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 254711.25
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 254711.25
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 2292401.25
$ws.Range("N72").Value = -2301137.25
$ws.Range("M72").ClearContents()

$ws.Range("H129").Value = 41482988
$ws.Range("J129").Value = 1853606.2
$ws.Range("L129").Value = 5560818.6
$ws.Range("N129").Value = -5570818.6

$ws.Range("H132").Value = 2203.4
$ws.Range("I132").Value = 2089.3572
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 6268.071599999999
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -3738.071599999999
$ws.Range("N132").Value = -16460

$ws.Range("H137").Value = 10205704
$ws.Range("I137").Value = 1550.8438
$ws.Range("J137").Value = 29413520
$ws.Range("K137").Value = 4652.5314
$ws.Range("L137").Value = 88240560
$ws.Range("M137").Value = -2102.5314
$ws.Range("N137").Value = -88245660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23524.318
$ws.Range("I32").Value = 20196.688
$ws.Range("K32").Value = 20196.688
$ws.Range("M32").Value = -19909.688

$ws.Range("H132").Value = 2374.6453
$ws.Range("I132").Value = 1766.8889
$ws.Range("J132").Value = 3216.1538
$ws.Range("K132").Value = 5300.6667
$ws.Range("L132").Value = 9648.4614
$ws.Range("M132").Value = -2770.6667
$ws.Range("N132").Value = -14708.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 37062012
$ws.Range("I20").Value = 59581.184
$ws.Range("J20").Value = 62501180
$ws.Range("K20").Value = 59581.184
$ws.Range("L20").Value = 62501180
$ws.Range("M20").Value = -59334.184
$ws.Range("N20").Value = -62501674

$ws.Range("H86").Value = 2247.963
$ws.Range("I86").Value = 2046.9474
$ws.Range("J86").Value = 2725.375
$ws.Range("K86").Value = 2046.9474
$ws.Range("L86").Value = 2725.375
$ws.Range("M86").Value = -923.9474
$ws.Range("N86").Value = -4971.375

$ws.Range("H89").Value = 2247.963
$ws.Range("I89").Value = 2046.9474
$ws.Range("J89").Value = 2725.375
$ws.Range("K89").Value = 10234.737
$ws.Range("L89").Value = 13626.875
$ws.Range("M89").Value = -4618.737000000001
$ws.Range("N89").Value = -24858.875

$ws.Range("H94").Value = 1317.579
$ws.Range("I94").Value = 1294
$ws.Range("J94").Value = 1383.6
$ws.Range("K94").Value = 1294
$ws.Range("L94").Value = 1383.6
$ws.Range("M94").Value = -843
$ws.Range("N94").Value = -2285.6

$ws.Range("H107").Value = 13802.637
$ws.Range("I107").Value = 1875.375
$ws.Range("J107").Value = 45608.668
$ws.Range("K107").Value = 1875.375
$ws.Range("L107").Value = 45608.668
$ws.Range("M107").Value = 44.625
$ws.Range("N107").Value = -49448.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 148.36842
$ws.Range("I7").Value = 128.46666
$ws.Range("J7").Value = 223
$ws.Range("K7").Value = 128.46666
$ws.Range("L7").Value = 223
$ws.Range("M7").Value = -15.46665999999999
$ws.Range("N7").Value = -449

$ws.Range("H31").Value = 1603.0209
$ws.Range("I31").Value = 1007.9474
$ws.Range("J31").Value = 1992.8966
$ws.Range("K31").Value = 1007.9474
$ws.Range("L31").Value = 1992.8966
$ws.Range("M31").Value = -712.9474
$ws.Range("N31").Value = -2582.8966

$ws.Range("H34").Value = 1603.0209
$ws.Range("I34").Value = 1007.9474
$ws.Range("J34").Value = 1992.8966
$ws.Range("K34").Value = 1007.9474
$ws.Range("L34").Value = 1992.8966
$ws.Range("M34").Value = -805.9474
$ws.Range("N34").Value = -2396.8966

$ws.Range("H140").Value = 48048
$ws.Range("J140").Value = 48048
$ws.Range("L140").Value = 48048
$ws.Range("N140").Value = -58408

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1349.5405
$ws.Range("I68").Value = 717.4194
$ws.Range("J68").Value = 1805.2559
$ws.Range("K68").Value = 2152.2582
$ws.Range("L68").Value = 5415.7677
$ws.Range("M68").Value = -1341.2582
$ws.Range("N68").Value = -7037.7677

$ws.Range("H71").Value = 1349.5405
$ws.Range("I71").Value = 717.4194
$ws.Range("J71").Value = 1805.2559
$ws.Range("K71").Value = 6456.7746
$ws.Range("L71").Value = 16247.3031
$ws.Range("M71").Value = -2400.7746
$ws.Range("N71").Value = -24359.3031

$ws.Range("H107").Value = 870.625
$ws.Range("J107").Value = 1616.8462
$ws.Range("L107").Value = 4850.5386
$ws.Range("N107").Value = -8690.5386

$ws.Range("H113").Value = 445.27835
$ws.Range("I113").Value = 392.64706
$ws.Range("J113").Value = 456.4625
$ws.Range("K113").Value = 1177.94118
$ws.Range("L113").Value = 1369.3875
$ws.Range("M113").Value = 992.05882
$ws.Range("N113").Value = -5709.3875

$ws.Range("H131").Value = 15762.112
$ws.Range("J131").Value = 1764.8448
$ws.Range("L131").Value = 5294.5344
$ws.Range("N131").Value = -15374.5344

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1672.8823
$ws.Range("I97").Value = 1713.5
$ws.Range("J97").Value = 1483.3334
$ws.Range("K97").Value = 1713.5
$ws.Range("L97").Value = 1483.3334
$ws.Range("M97").Value = -1217.5
$ws.Range("N97").Value = -2475.3334

$ws.Range("H113").Value = 1527.3
$ws.Range("I113").Value = 1183.3334
$ws.Range("J113").Value = 2043.25
$ws.Range("K113").Value = 1183.3334
$ws.Range("L113").Value = 2043.25
$ws.Range("M113").Value = 986.6666
$ws.Range("N113").Value = -6383.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 336666.66
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5386

$ws.Range("H61").Value = 909.4545000000001
$ws.Range("I61").Value = 888
$ws.Range("J61").Value = 966.6667
$ws.Range("K61").Value = 888
$ws.Range("L61").Value = 966.6667
$ws.Range("M61").Value = -686
$ws.Range("N61").Value = -1370.6667

$ws.Range("H113").Value = 909.4545000000001
$ws.Range("I113").Value = 888
$ws.Range("J113").Value = 966.6667
$ws.Range("K113").Value = 888
$ws.Range("L113").Value = 966.6667
$ws.Range("M113").Value = 1282
$ws.Range("N113").Value = -5306.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1281.8334
$ws.Range("I81").Value = 1284.75
$ws.Range("J81").Value = 1276
$ws.Range("K81").Value = 2569.5
$ws.Range("L81").Value = 2552
$ws.Range("M81").Value = -1508.5
$ws.Range("N81").Value = -4674

$ws.Range("H84").Value = 1281.8334
$ws.Range("I84").Value = 1284.75
$ws.Range("J84").Value = 1276
$ws.Range("K84").Value = 12847.5
$ws.Range("L84").Value = 12760
$ws.Range("M84").Value = -7543.5
$ws.Range("N84").Value = -23368
